$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.32 = 50363.56 pesos`n✅ 50363.56 pesos = 12.4 = 984.94 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the "tasas" sheet: updated rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 81.17
$wsTasas.Range("O10").Value = 4088.01
$wsTasas.Range("N12").Value = 4060
$wsTasas.Range("O12").Value = 79.40000000000001
